# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# "Price" (column D) cells are stored as text in this sheet, even when the
# text looks numeric (e.g. "217.31", or "0.07800" with a meaningful
# trailing zero) so a leading apostrophe is used to force literal text
# instead of having Excel reinterpret/round it as a number; the cell
# style is then reset to "Normal" so the apostrophe's "number stored as
# text" indicator doesn't leave a lingering style change behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.096.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "'1.653.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'217.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'0.5255"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.2597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "'0.07800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "'4.498"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "'1.586.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("D14").Value = "'0.5476"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "'0.0₅8226"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "'65.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "'26.098.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'4.574"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").Value = "'190.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'6.023"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'141.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "'0.1232"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "'7.235"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'16.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "'1.429"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "'0.05837"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").Value = "'1.273"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'3.530"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").Value = "'3.258"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "'0.9479"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'2.413"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'2.778"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'0.5729"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'0.01610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "'5.776"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.8450"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'103.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("D43").Value = "'1.025.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "'1.798.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'57.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").Value = "'0.4313"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'0.05148"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.467"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.786"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "'0.09657"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.51%  "
